$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3399.1333
$ws.Cells.Item(76, 9).Value = 2500.3333
$ws.Cells.Item(76, 11).Value = 2500.3333
$ws.Cells.Item(76, 13).Value = -2185.3333

$ws.Cells.Item(79, 8).Value = 3399.1333
$ws.Cells.Item(79, 9).Value = 2500.3333
$ws.Cells.Item(79, 11).Value = 2500.3333
$ws.Cells.Item(79, 13).Value = -1408.3333

$ws.Cells.Item(138, 8).Value = 3516.6572
$ws.Cells.Item(138, 9).Value = 1229.5625
$ws.Cells.Item(138, 10).Value = 4194.315
$ws.Cells.Item(138, 11).Value = 3688.6875
$ws.Cells.Item(138, 12).Value = 12582.945
$ws.Cells.Item(138, 13).Value = 1451.3125
$ws.Cells.Item(138, 14).Value = -22862.945

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1538.102
$ws.Cells.Item(20, 9).Value = 1090
$ws.Cells.Item(20, 10).Value = 2309.8333
$ws.Cells.Item(20, 11).Value = 1090
$ws.Cells.Item(20, 12).Value = 2309.8333
$ws.Cells.Item(20, 13).Value = -843
$ws.Cells.Item(20, 14).Value = -2803.8333

$ws.Cells.Item(99, 8).Value = 2599
$ws.Cells.Item(99, 9).Value = 2218.9
$ws.Cells.Item(99, 10).Value = 3142
$ws.Cells.Item(99, 11).Value = 2218.9
$ws.Cells.Item(99, 12).Value = 3142
$ws.Cells.Item(99, 13).Value = -720.9000000000001
$ws.Cells.Item(99, 14).Value = -6138

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2937.5557
$ws.Cells.Item(16, 9).Value = 1670.6666
$ws.Cells.Item(16, 10).Value = 3571
$ws.Cells.Item(16, 11).Value = 1670.6666
$ws.Cells.Item(16, 12).Value = 3571
$ws.Cells.Item(16, 13).Value = -1383.6666
$ws.Cells.Item(16, 14).Value = -4145

$ws.Cells.Item(31, 8).Value = 2258.509
$ws.Cells.Item(31, 9).Value = 1517.6552
$ws.Cells.Item(31, 10).Value = 3084.8462
$ws.Cells.Item(31, 11).Value = 1517.6552
$ws.Cells.Item(31, 12).Value = 3084.8462
$ws.Cells.Item(31, 13).Value = -1222.6552
$ws.Cells.Item(31, 14).Value = -3674.8462

$ws.Cells.Item(34, 8).Value = 2258.509
$ws.Cells.Item(34, 9).Value = 1517.6552
$ws.Cells.Item(34, 10).Value = 3084.8462
$ws.Cells.Item(34, 11).Value = 1517.6552
$ws.Cells.Item(34, 12).Value = 3084.8462
$ws.Cells.Item(34, 13).Value = -1315.6552
$ws.Cells.Item(34, 14).Value = -3488.8462

$ws.Cells.Item(58, 8).Value = 1042.3823
$ws.Cells.Item(58, 9).Value = 620.05554
$ws.Cells.Item(58, 10).Value = 2671.3572
$ws.Cells.Item(58, 11).Value = 620.05554
$ws.Cells.Item(58, 12).Value = 2671.3572
$ws.Cells.Item(58, 13).Value = -417.05554
$ws.Cells.Item(58, 14).Value = -3077.3572

$ws.Cells.Item(86, 8).Value = 5435.9165
$ws.Cells.Item(86, 10).Value = 3390.5
$ws.Cells.Item(86, 12).Value = 3390.5
$ws.Cells.Item(86, 14).Value = -5636.5

$ws.Cells.Item(89, 8).Value = 5435.9165
$ws.Cells.Item(89, 10).Value = 3390.5
$ws.Cells.Item(89, 12).Value = 16952.5
$ws.Cells.Item(89, 14).Value = -28184.5

$ws.Cells.Item(107, 8).Value = 778.5909
$ws.Cells.Item(107, 9).Value = 490.3
$ws.Cells.Item(107, 10).Value = 1018.8333
$ws.Cells.Item(107, 11).Value = 490.3
$ws.Cells.Item(107, 12).Value = 1018.8333
$ws.Cells.Item(107, 13).Value = 1429.7
$ws.Cells.Item(107, 14).Value = -4858.8333

$ws.Cells.Item(113, 8).Value = 2937.5557
$ws.Cells.Item(113, 9).Value = 1670.6666
$ws.Cells.Item(113, 10).Value = 3571
$ws.Cells.Item(113, 11).Value = 1670.6666
$ws.Cells.Item(113, 12).Value = 3571
$ws.Cells.Item(113, 13).Value = 499.3334
$ws.Cells.Item(113, 14).Value = -7911

$ws.Cells.Item(122, 8).Value = 1969.3334
$ws.Cells.Item(122, 9).Value = 1237.3334
$ws.Cells.Item(122, 10).Value = 3433.3333
$ws.Cells.Item(122, 11).Value = 3712.0002
$ws.Cells.Item(122, 12).Value = 10299.9999
$ws.Cells.Item(122, 13).Value = -1262.0002
$ws.Cells.Item(122, 14).Value = -15199.9999

$ws.Cells.Item(132, 8).Value = 2219.524
$ws.Cells.Item(132, 9).Value = 1567.7646
$ws.Cells.Item(132, 10).Value = 4989.5
$ws.Cells.Item(132, 11).Value = 4703.293799999999
$ws.Cells.Item(132, 12).Value = 14968.5
$ws.Cells.Item(132, 13).Value = -2173.293799999999
$ws.Cells.Item(132, 14).Value = -20028.5

$ws.Cells.Item(134, 8).Value = 2228.1365
$ws.Cells.Item(134, 9).Value = 1129.7059
$ws.Cells.Item(134, 10).Value = 5962.8
$ws.Cells.Item(134, 11).Value = 3389.1177
$ws.Cells.Item(134, 12).Value = 17888.4
$ws.Cells.Item(134, 13).Value = -854.1176999999998
$ws.Cells.Item(134, 14).Value = -22958.4

$ws.Cells.Item(136, 8).Value = 1042.3823
$ws.Cells.Item(136, 9).Value = 620.05554
$ws.Cells.Item(136, 10).Value = 2671.3572
$ws.Cells.Item(136, 11).Value = 1860.16662
$ws.Cells.Item(136, 12).Value = 8014.071599999999
$ws.Cells.Item(136, 13).Value = 689.83338
$ws.Cells.Item(136, 14).Value = -13114.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 418.18182
$ws.Cells.Item(68, 9).Value = 375
$ws.Cells.Item(68, 10).Value = 442.85715
$ws.Cells.Item(68, 11).Value = 1125
$ws.Cells.Item(68, 12).Value = 1328.57145
$ws.Cells.Item(68, 13).Value = -314
$ws.Cells.Item(68, 14).Value = -2950.57145

$ws.Cells.Item(71, 8).Value = 418.18182
$ws.Cells.Item(71, 9).Value = 375
$ws.Cells.Item(71, 10).Value = 442.85715
$ws.Cells.Item(71, 11).Value = 3375
$ws.Cells.Item(71, 12).Value = 3985.71435
$ws.Cells.Item(71, 13).Value = 681
$ws.Cells.Item(71, 14).Value = -12097.71435

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

$ws.Cells.Item(86, 8).Value = 1059.4736
$ws.Cells.Item(86, 9).Value = 803.3333
$ws.Cells.Item(86, 10).Value = 1290
$ws.Cells.Item(86, 11).Value = 2409.9999
$ws.Cells.Item(86, 12).Value = 3870
$ws.Cells.Item(86, 13).Value = -1223.9999
$ws.Cells.Item(86, 14).Value = -6242

$ws.Cells.Item(89, 8).Value = 1059.4736
$ws.Cells.Item(89, 9).Value = 803.3333
$ws.Cells.Item(89, 10).Value = 1290
$ws.Cells.Item(89, 11).Value = 7229.9997
$ws.Cells.Item(89, 12).Value = 11610
$ws.Cells.Item(89, 13).Value = -1301.9997
$ws.Cells.Item(89, 14).Value = -23466

$ws.Cells.Item(92, 8).Value = 763
$ws.Cells.Item(92, 10).Value = 763
$ws.Cells.Item(92, 12).Value = 2289
$ws.Cells.Item(92, 14).Value = -4785

$ws.Cells.Item(107, 8).Value = 325
$ws.Cells.Item(107, 9).Value = 200
$ws.Cells.Item(107, 10).Value = 700
$ws.Cells.Item(107, 11).Value = 600
$ws.Cells.Item(107, 12).Value = 2100
$ws.Cells.Item(107, 13).Value = 1320
$ws.Cells.Item(107, 14).Value = -5940

$ws.Cells.Item(122, 8).Value = 12024.211
$ws.Cells.Item(122, 9).Value = 14011.625
$ws.Cells.Item(122, 10).Value = 1424.6666
$ws.Cells.Item(122, 11).Value = 126104.625
$ws.Cells.Item(122, 12).Value = 12821.9994
$ws.Cells.Item(122, 13).Value = -123654.625
$ws.Cells.Item(122, 14).Value = -17721.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5579.857
$ws.Cells.Item(70, 9).Value = 5626.5
$ws.Cells.Item(70, 11).Value = 5626.5
$ws.Cells.Item(70, 13).Value = -5356.5

$ws.Cells.Item(73, 8).Value = 5579.857
$ws.Cells.Item(73, 9).Value = 5626.5
$ws.Cells.Item(73, 11).Value = 5626.5
$ws.Cells.Item(73, 13).Value = -4690.5

$ws.Cells.Item(113, 8).Value = 8558.866
$ws.Cells.Item(113, 9).Value = 1688.5
$ws.Cells.Item(113, 10).Value = 13139.111
$ws.Cells.Item(113, 11).Value = 1688.5
$ws.Cells.Item(113, 12).Value = 13139.111
$ws.Cells.Item(113, 13).Value = 481.5
$ws.Cells.Item(113, 14).Value = -17479.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1378.3158
$ws.Cells.Item(16, 9).Value = 511.75
$ws.Cells.Item(16, 10).Value = 6000
$ws.Cells.Item(16, 11).Value = 511.75
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = -341.75
$ws.Cells.Item(16, 14).Value = -6340

$ws.Cells.Item(22, 8).Value = 678.17645
$ws.Cells.Item(22, 9).Value = 684.4545000000001
$ws.Cells.Item(22, 10).Value = 666.6667
$ws.Cells.Item(22, 11).Value = 684.4545000000001
$ws.Cells.Item(22, 12).Value = 666.6667
$ws.Cells.Item(22, 13).Value = -389.4545000000001
$ws.Cells.Item(22, 14).Value = -1256.6667

$ws.Cells.Item(27, 8).Value = 678.17645
$ws.Cells.Item(27, 9).Value = 684.4545000000001
$ws.Cells.Item(27, 10).Value = 666.6667
$ws.Cells.Item(27, 11).Value = 684.4545000000001
$ws.Cells.Item(27, 12).Value = 666.6667
$ws.Cells.Item(27, 13).Value = -577.4545000000001
$ws.Cells.Item(27, 14).Value = -880.6667

$ws.Cells.Item(46, 8).Value = 111111890
$ws.Cells.Item(46, 9).Value = 786
$ws.Cells.Item(46, 10).Value = 500000740
$ws.Cells.Item(46, 11).Value = 786
$ws.Cells.Item(46, 12).Value = 500000740
$ws.Cells.Item(46, 13).Value = -598
$ws.Cells.Item(46, 14).Value = -500001116

$ws.Cells.Item(55, 8).Value = 17544496
$ws.Cells.Item(55, 9).Value = 23810172
$ws.Cells.Item(55, 10).Value = 605
$ws.Cells.Item(55, 11).Value = 23810172
$ws.Cells.Item(55, 12).Value = 605
$ws.Cells.Item(55, 13).Value = -23809999
$ws.Cells.Item(55, 14).Value = -951

$ws.Cells.Item(61, 8).Value = 2035.3
$ws.Cells.Item(61, 9).Value = 1359.5
$ws.Cells.Item(61, 10).Value = 2485.8333
$ws.Cells.Item(61, 11).Value = 1359.5
$ws.Cells.Item(61, 12).Value = 2485.8333
$ws.Cells.Item(61, 13).Value = -1157.5
$ws.Cells.Item(61, 14).Value = -2889.8333

$ws.Cells.Item(113, 8).Value = 2035.3
$ws.Cells.Item(113, 9).Value = 1359.5
$ws.Cells.Item(113, 10).Value = 2485.8333
$ws.Cells.Item(113, 11).Value = 1359.5
$ws.Cells.Item(113, 12).Value = 2485.8333
$ws.Cells.Item(113, 13).Value = 810.5
$ws.Cells.Item(113, 14).Value = -6825.8333

$ws.Cells.Item(122, 8).Value = 4045
$ws.Cells.Item(122, 9).Value = 4100
$ws.Cells.Item(122, 10).Value = 3990
$ws.Cells.Item(122, 11).Value = 12300
$ws.Cells.Item(122, 12).Value = 11970
$ws.Cells.Item(122, 13).Value = -9850
$ws.Cells.Item(122, 14).Value = -16870

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1026.1187
$ws.Cells.Item(132, 9).Value = 608.1
$ws.Cells.Item(132, 10).Value = 1906.1578
$ws.Cells.Item(132, 11).Value = 1824.3
$ws.Cells.Item(132, 12).Value = 5718.4734
$ws.Cells.Item(132, 13).Value = 705.6999999999998
$ws.Cells.Item(132, 14).Value = -10778.4734
